# Regenerate test data to include credit card account refund/credit
# transactions: the Expense/Income/Transfer pivot sheets are refreshed
# with new per-period totals (and their Total row/column) reflecting the
# added refund/credit transactions on credit card accounts.

$wb = $excel.ActiveWorkbook

# --- Expense sheet (rows 2-9 = categories, row 10 = Total; col G = Total) ---
$ws = $wb.Worksheets.Item("Expense")
$expenseData = @{
    "B2" = -38915.87;  "C2" = -42885.77;  "D2" = -38293.06;  "E2" = -31954.07;  "F2" = -28973.95;  "G2" = -181022.72
    "B3" = -36353.05;  "C3" = -30457.38;  "D3" = -38819.41;  "E3" = -38946.09;  "F3" = -34430.07;  "G3" = -179006
    "B4" = -11866.61;  "C4" = -20344.57;  "D4" = -15438.74;  "E4" = -15008.15;  "F4" = -18061.25;  "G4" = -80719.32000000001
    "B5" = -42160.42;  "C5" = -40308.66;  "D5" = -31167.05;  "E5" = -37792.51;  "F5" = -39080.4;   "G5" = -190509.04
    "B6" = -11891.18;  "C6" = -20277.56;  "D6" = -16132.08;  "E6" = -18620.19;  "F6" = -13311.51;  "G6" = -80232.52
    "B7" = -57936.9;   "C7" = -44474.84;  "D7" = -54028.1;   "E7" = -73444.08;  "F7" = -49862.93;  "G7" = -279746.85
    "B8" = -28564.63;  "C8" = -29994.23;  "D8" = -27644.43;  "E8" = -29930.08;  "F8" = -24965.85;  "G8" = -141099.22
    "B9" = -18135.67;  "C9" = -16482.58;  "D9" = -15629.56;  "E9" = -19907.67;  "F9" = -12846.96;  "G9" = -83002.44
    "B10" = -245824.33; "C10" = -245225.59; "D10" = -237152.43; "E10" = -265602.84; "F10" = -221532.92; "G10" = -1215338.11
}
foreach ($addr in $expenseData.Keys) {
    $ws.Range($addr).Value = $expenseData[$addr]
}

# --- Income sheet (rows 2-4 = categories, row 5 = Total; col G = Total) ---
$ws = $wb.Worksheets.Item("Income")
$incomeData = @{
    "B2" = 57327.05;  "C2" = 69089.91;  "D2" = 46505.97;  "E2" = 68042.99000000001; "F2" = 97647.47;  "G2" = 338613.39
    "B3" = 452802.19; "C3" = 439948.04; "D3" = 521104.54; "E3" = 400785.61;         "F3" = 439761.29; "G3" = 2254401.67
    "B4" = 81716.10000000001; "C4" = 76255.21000000001; "D4" = 76187.89999999999; "E4" = 68580.99000000001; "F4" = 110584.34; "G4" = 413324.54
    "B5" = 591845.34; "C5" = 585293.16; "D5" = 643798.41; "E5" = 537409.59;        "F5" = 647993.1;  "G5" = 3006339.6
}
foreach ($addr in $incomeData.Keys) {
    $ws.Range($addr).Value = $incomeData[$addr]
}

# --- Transfer sheet (rows 2-3 = Transfer From/To; col G = Total) ---
$ws = $wb.Worksheets.Item("Transfer")
$transferData = @{
    "B2" = 196582.08;  "C2" = 188143.44;  "D2" = 185434.74;  "E2" = 205682.74;  "F2" = 189920.85;  "G2" = 965763.85
    "B3" = -196582.08; "C3" = -188143.44; "D3" = -185434.74; "E3" = -205682.74; "F3" = -189920.85; "G3" = -965763.85
}
foreach ($addr in $transferData.Keys) {
    $ws.Range($addr).Value = $transferData[$addr]
}
